# Lowercase every "subcategory" value in column H (rows 2-65) of Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 8).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 8)
    $value = $cell.Value()
    if ($null -ne $value -and $value -is [string]) {
        $cell.Value = $value.ToLower()
    }
}
